$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the environment data (B2:E2, H2) to the new Preprod environment values.
# B2 carries a quote-prefix style (s="1"); prefixing the new value with an
# apostrophe preserves that text-entry style the way typing it in Excel would.
$ws.Range("B2").Value = "'i-preproducciongestion.segurossura.com.ar"
$ws.Range("C2").Value = "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"
$ws.Range("E2").Value = "silverarrow"
$ws.Range("H2").Value = 2240451788

# Remove the hyperlink on C2 (keep the text/value, drop the link)
$ws.Hyperlinks.Delete()

# Update the active selection to W2
$ws.Range("W2").Select()
